# Updates cryptos list prices/volumes to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.175.69"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.641.24"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D5").Value = "217.01"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D9").Value = "0.0626"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D10").Value = "19.95"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.872.20"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "1.642.03"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D15").Value = "0.541"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D16").Value = "67.04"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "27.183.98"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D19").Value = "217.94"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D21").Value = "6.95"
$ws.Range("E21").Value = "  +2.76%  "
$ws.Range("D22").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D22").Value = "4.41"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("D24").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D25").Value = "147.22"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D27").Value = "7.44"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D28").Value = "0.119"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D29").Value = "15.70"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D30").Value = "0.0508"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D34").Value = "1.57"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "1.299.65"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D38").Value = "0.548"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D39").Value = "0.856"
$ws.Range("E39").Value = "  +2.99%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +6.03%  "
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").Value = "1.782.29"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D45").Value = "61.71"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D46").Value = "91.78"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D50").Value = "7.67"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"  # keep as text, not auto-converted to a number
$ws.Range("D51").Value = "0.0963"
$ws.Range("E51").Value = "  +0.30%  "
